# Apply the invoice edit described by the commit:
#   "Added table, header, total price, company and logo."
#
# Concretely (per the xlsx XML diff) the meaningful content changes are:
#   1. Header cell C1 is retitled from "amount_purchased" to "purchased".
#   2. A new line-item row (row 4) is added to the table:
#        product_id=4567757, product_name="PS5 Pro",
#        purchased=2, price_per_unit=800, total_price=1600
#      formatted like the preceding data rows (bold/text id cell in column A,
#      matching the row above; plain cells elsewhere).
#   3. Column C is widened (auto-fit-ish) to comfortably fit the new header.
#   4. The active selection ends up on C1 (where the header was edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the "amount_purchased" header to "purchased".
$ws.Range("C1").Value = "purchased"

# 2. Add the new product row (row 4) under the existing two line items.
#    Copy the formatting of the row above (row 3) for the product-id cell so
#    the bold/text styling carries down, matching the existing table rows.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A4").Value = "4567757"
$ws.Range("B4").Value = "PS5 Pro"
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 800
$ws.Range("E4").Value = 1600

# 3. Widen column C a bit so the new "purchased" header fits comfortably
#    (mirrors the best-fit column resize visible in the saved file).
$ws.Columns.Item(3).ColumnWidth = 120.5 / 7

# 4. Leave the selection on the header cell that was just edited.
[void]$ws.Range("C1").Select()
